$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full data block (rows 2-9), columns A:T
# Columns: A=Sending cluster, B=Ligand symbol, C=Receptor symbol, D=Target cluster,
# E..T = numeric metrics
$data = @(
    @("ECs",  "Ccl11", "Cxcr3", "FAPs", 2, 0.6666666666666666, 3.498267, 10.494801, 0.02652289550929144, 0.02652289550929144, 1, 0.3333333333333333, 0.05194133333333333, 0.155824, 0.03788844568234288, 0.03788844568234288, 0.181704652336, 1.635341871024, 0.001004911285842245, 0.001004911285842245),
    @("ECs",  "Ccl11", "Cxcr3", "M2",   2, 0.6666666666666666, 3.498267, 10.494801, 0.02652289550929144, 0.02652289550929144, 3, 1, 1.318960333333333, 3.956881, 0.962111554317657, 0.9621115543176572, 4.614075408409, 41.526678675681, 0.02551798422344919, 0.0255179842234492),
    @("FAPs", "Ccl11", "Cxcr3", "FAPs", 3, 1, 122.458089, 367.374267, 0.9284434546632696, 0.9284434546632696, 1, 0.3333333333333333, 0.05194133333333333, 0.155824, 0.03788844568234288, 0.03788844568234288, 6.360636420112001, 57.245727781008, 0.03517727940113607, 0.03517727940113607),
    @("FAPs", "Ccl11", "Cxcr3", "M2",   3, 1, 122.458089, 367.374267, 0.9284434546632696, 0.9284434546632696, 3, 1, 1.318960333333333, 3.956881, 0.962111554317657, 0.9621115543176572, 161.517361886803, 1453.656256981227, 0.8932661752621335, 0.8932661752621336),
    @("M2",   "Ccl11", "Cxcr3", "FAPs", 3, 1, 0.474194, 1.422582, 0.003595208116799817, 0.003595208116799817, 1, 0.3333333333333333, 0.05194133333333333, 0.155824, 0.03788844568234288, 0.03788844568234288, 0.02463026861866667, 0.221672417568, 0.0001362168474500881, 0.0001362168474500881),
    @("M2",   "Ccl11", "Cxcr3", "M2",   3, 1, 0.474194, 1.422582, 0.003595208116799817, 0.003595208116799817, 3, 1, 1.318960333333333, 3.956881, 0.962111554317657, 0.9621115543176572, 0.6254430763046667, 5.628987686742001, 0.003458991269349729, 0.003458991269349729),
    @("sCs",  "Ccl11", "Cxcr3", "FAPs", 3, 1, 5.465569666666667, 16.396709, 0.04143844171063926, 0.04143844171063926, 1, 0.3333333333333333, 0.05194133333333333, 0.155824, 0.03788844568234288, 0.03788844568234288, 0.2838889759128889, 2.555000783216, 0.001570038147914487, 0.001570038147914487),
    @("sCs",  "Ccl11", "Cxcr3", "M2",   3, 1, 5.465569666666667, 16.396709, 0.04143844171063926, 0.04143844171063926, 3, 1, 1.318960333333333, 3.956881, 0.962111554317657, 0.9621115543176572, 7.208869589403222, 64.87982630462901, 0.03986840356272477, 0.03986840356272477)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    for ($c = 0; $c -lt $vals.Count; $c++) {
        $col = $c + 1
        $ws.Cells.Item($row, $col).Value = $vals[$c]
    }
}
